$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "empty but formatted" cells (style index 0 - the plain/general style
#    already used by cells such as B4, G4, H4, I4, ...). These cells are
#    being added to columns J/K (Trainer / Customer) on rows 3-12 to line up
#    the grid with the two new class columns, without marking an "X".
# ---------------------------------------------------------------------------
$s0src = $ws.Range("B4")

$s0src.Copy()
$ws.Range("J3").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J4:K4").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("K5").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J6").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J7:K7").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J8:K8").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J9:K9").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J10").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J11").PasteSpecial(-4122)

$s0src.Copy()
$ws.Range("J12:K12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Row 14 (UC-13): the "X" mark moves from column B (SecurityController)
#    to column C (UserController). Clear B14's value but keep its format,
#    then write the "X" into C14 and restore the row's bold/centered style.
# ---------------------------------------------------------------------------
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "X"
$ws.Range("E14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) New use cases appended as rows 15-17.
# ---------------------------------------------------------------------------

# Style helpers pulled from existing rows:
#   A14  -> s=1 (left aligned) row-label style, used for the new A15:A17
#   D13  -> s=3 with no value (blank bold/centered cell)
#   E14  -> s=3 with t="s" v="21" ("X" mark, bold/centered)

# --- Row 15: UC-14 - Modify User Account ---
$ws.Range("A15").Value = "UC-14 – Modify User Account"
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("D13").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C15").Value = "X"
$ws.Range("E15:K15").Value = "X"
$ws.Range("E14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15:K15").PasteSpecial(-4122)

# --- Row 16: UC-15 - Search for Trainers ---
$ws.Range("A16").Value = "UC-15 – Search for Trainers"
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("C16").Value = "X"
$ws.Range("E16").Value = "X"
$ws.Range("G16").Value = "X"
$ws.Range("J16").Value = "X"
$ws.Range("E14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("J16").PasteSpecial(-4122)

$ws.Range("D13").Copy()
$ws.Range("F16").PasteSpecial(-4122)

# --- Row 17: UC-16 - Search for Equipment Items ---
$ws.Range("A17").Value = "UC-16 – Search for Equipment Items"
$ws.Range("A14").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("D17").Value = "X"
$ws.Range("E17").Value = "X"
$ws.Range("O17").Value = "X"
$ws.Range("E14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("O17").PasteSpecial(-4122)

$ws.Range("D13").Copy()
$ws.Range("F17").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Refresh the view: active cell/selection moves to U17, matching the
#    author's final cursor position after adding the new rows.
# ---------------------------------------------------------------------------
$ws.Range("U17").Select()
